# "Otp under 5 minutes" — recolor the "Reset Password with OTP" list items
# that are currently black (theme color "Text 1") to purple (#7030A0),
# matching the sibling "Reset Password with OTP" items that already use
# an explicit accent color. Leaves the already-blue (#0070C0) items and
# any unrelated black text untouched.

$d = $word.ActiveDocument

$newColor = 10498160  # RGB(0x70,0x30,0xA0) == #7030A0 as a VBA/wdColor long

foreach ($p in $d.Paragraphs) {
    $range = $p.Range
    if ($range.Text -like "*Reset Password with OTP*") {
        if ($range.Font.Color -eq -587137025) {
            $range.Font.Color = $newColor
        }
    }
}
